$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 310; existing rows 310:432 shift down to 311:433.
$ws.Rows(310).Insert()

# Populate the newly inserted row 310 with the new week's data.
$ws.Cells.Item(310, 1).Value = 10
$ws.Cells.Item(310, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(310, 3).Value = "La Araucanía"
$ws.Cells.Item(310, 4).Value = 44924
$ws.Cells.Item(310, 5).Value = 9
$ws.Cells.Item(310, 6).Value = 100112009
$ws.Cells.Item(310, 7).Value = "Acelga"
$ws.Cells.Item(310, 8).Value = "Sin especificar"
$ws.Cells.Item(310, 9).Value = "Primera"
$ws.Cells.Item(310, 10).Value = 115
$ws.Cells.Item(310, 11).Value = 9000
$ws.Cells.Item(310, 12).Value = 10000
$ws.Cells.Item(310, 13).Value = 9565
$ws.Cells.Item(310, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(310, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(310, 16).Value = 797
$ws.Cells.Item(310, 17).Value = 12
$ws.Cells.Item(310, 18).Value = "Hortaliza"

# Match the date formatting used by the rest of column D.
$ws.Cells.Item(310, 4).NumberFormat = $ws.Cells.Item(311, 4).NumberFormat
